$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# New row 166 (NUMBER 165)
$ws.Cells.Item(166, 1).Value = 165.0
$ws.Cells.Item(166, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(166, 3).Value = "5:55 PM"
$ws.Cells.Item(166, 4).Value = "FR1056"
$ws.Cells.Item(166, 5).Value = "Brussels"
$ws.Cells.Item(166, 6).Value = "(CRL)"
$ws.Cells.Item(166, 7).Value = "Ryanair "
$ws.Cells.Item(166, 8).Value = "B38M"
$ws.Cells.Item(166, 9).Value = "(SP-RZO)"
$ws.Cells.Item(166, 10).Value = "6:05 PM"
$ws.Cells.Item(166, 12).Value = "0 hours, 10 minutes"

# New row 167 (NUMBER 166)
$ws.Cells.Item(167, 1).Value = 166.0
$ws.Cells.Item(167, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(167, 3).Value = "6:20 PM"
$ws.Cells.Item(167, 4).Value = "FR1968"
$ws.Cells.Item(167, 5).Value = "Madrid"
$ws.Cells.Item(167, 6).Value = "(MAD)"
$ws.Cells.Item(167, 7).Value = "Ryanair "
$ws.Cells.Item(167, 8).Value = "B738"
$ws.Cells.Item(167, 9).Value = "(EI-ESV)"
$ws.Cells.Item(167, 10).Value = "6:34 PM"
$ws.Cells.Item(167, 12).Value = "0 hours, 14 minutes"

# New row 168 (NUMBER 167)
$ws.Cells.Item(168, 1).Value = 167.0
$ws.Cells.Item(168, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(168, 3).Value = "7:15 PM"
$ws.Cells.Item(168, 4).Value = "FR1922"
$ws.Cells.Item(168, 5).Value = "Budapest"
$ws.Cells.Item(168, 6).Value = "(BUD)"
$ws.Cells.Item(168, 7).Value = "Ryanair "
$ws.Cells.Item(168, 8).Value = "B738"
$ws.Cells.Item(168, 9).Value = "(SP-RSS)"
$ws.Cells.Item(168, 10).Value = "7:10 PM"
$ws.Cells.Item(168, 12).Value = "0 hours, -5 minutes"
